# Restore C10 ("From" value for rule R20) from 18 to 1, matching the
# committed revision (admin SAVE on 07/07/2020).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
